$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the "AUTos keret" column
$ws.Range("E2").Value = "AUTos keret"

# New rows of spending entries
$ws.Range("A4").Value = "Vonalszenzor kártya"
$ws.Range("E4").Value = 30000

$ws.Range("A5").Value = "Motormeghajtó kártya"
$ws.Range("E5").Value = 12000

# Match header styling (fill) used by the other header cells in row 2
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial(-4122)

# Column E width tweak, matching the recorded custom width
$ws.Columns.Item(5).ColumnWidth = 11.1

# Update the active selection, as captured in the saved view state
$ws.Range("F5").Select()
